# Update the 'Instructions' sheet:
#  - bump the version string
#  - reword the "no headers" instruction
#  - add a new instruction row about not leaving blank rows, pushing the
#    remaining rows down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# The sheet is protected; unprotect so we can edit, then restore protection.
$ws.Unprotect()

$ws.Range("A2").Value = "Version 1.2.3"
$ws.Range("A5").Value = "Please use consecutive rows (no blank rows)."

# Insert a brand-new row 6; everything from the old row 6 onward shifts
# down by one (old row 6 "Do not edit the other sheets." becomes row 7, the
# old header block starting at row 8 becomes row 9, etc.)
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Do not edit the header row of the 'Antibodies' sheet."

$ws.Protect()
